# mod (Point) : ajout de operator*
#
# The class-diagram table cell for "Point" had every member line spread
# across several runs (with spell/grammar-check proofErr markers from the
# original authoring session). This edit retypes each line as clean,
# single-run text -- and, in the process, drops the leading "-" from the
# "operator-" line and inserts a brand-new "+operator*  (Point) : Point"
# line right after it.

$d = $word.ActiveDocument

function Replace-Line($oldText, $newText) {
    $range = $d.Content
    $ok = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $newText, 1)
    if (-not $ok) {
        Write-Output "NOT FOUND: $oldText"
    }
}

# Lines whose wording is unchanged but whose runs get consolidated.
Replace-Line "<<create>>+Point(double,double)" "<<create>>+Point(double,double)"
Replace-Line "+Mov(double,double) : void" "+Mov(double,double) : void"
Replace-Line "+Distancer(Point) : double" "+Distancer(Point) : double"
Replace-Line "+Add (Point) : void" "+Add (Point) : void"
Replace-Line "+operator+ (Point) : Point" "+operator+ (Point) : Point"
Replace-Line "+Affichage() : void" "+Affichage() : void"
Replace-Line "+GetX() : double" "+GetX() : double"
Replace-Line "+GetY() : double" "+GetY() : double"
Replace-Line "+SetX(double) : void" "+SetX(double) : void"
Replace-Line "+SetY(double) : void" "+SetY(double) : void"

# "-operator- (Point) : Point" loses its leading "-" and gains a new
# sibling paragraph for "+operator*  (Point) : Point" right after it.
# Rewriting this (single-run) paragraph's Range.Text with an embedded
# carriage return splits it into two paragraphs, and the newly created
# second paragraph correctly inherits the paragraph/run formatting of the
# paragraph mark it was split from.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($txt -eq "-operator- (Point) : Point") {
        $p.Range.Text = "operator- (Point) : Point" + [char]13 + "+operator*  (Point) : Point"
        break
    }
}
